$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 146, shifting existing rows 146-210 down to 147-211
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new record
$ws.Cells.Item(146, 1).Value = 5
$ws.Cells.Item(146, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(146, 3).Value = "Maule"
$ws.Cells.Item(146, 4).Value = 44636
$ws.Cells.Item(146, 5).Value = 7
$ws.Cells.Item(146, 6).Value = 100112009
$ws.Cells.Item(146, 7).Value = "Acelga"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 200
$ws.Cells.Item(146, 11).Value = 4000
$ws.Cells.Item(146, 12).Value = 4000
$ws.Cells.Item(146, 13).Value = 4000
$ws.Cells.Item(146, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(146, 15).Value = "Región del Maule"
$ws.Cells.Item(146, 16).Value = 1000
$ws.Cells.Item(146, 17).Value = 4
$ws.Cells.Item(146, 18).Value = "Hortaliza"

# Match the date style used by column D in the surrounding rows
$ws.Cells.Item(146, 4).NumberFormat = $ws.Cells.Item(145, 4).NumberFormat
